$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 5373.000859812833
$ws.Range("E4").Value = 3268.002014659345
$ws.Range("F4").Value = 0.6082265943976016
$ws.Range("G4").Value = 1.644124096530862
$ws.Range("H4").Value = 145.9361429969722
$ws.Range("I4").Value = 15.89677538211254
$ws.Range("J4").Value = 12.25906828093957
$ws.Range("K4").Value = 290.2080667279661
$ws.Range("L4").Value = 280.3078072129283
$ws.Range("M4").Value = 163.4015853721648
$ws.Range("N4").Value = 2.784134542991524
$ws.Range("O4").Value = 4.935462251400168
$ws.Range("P4").Value = 0.007548424894594064
$ws.Range("Q4").Value = 348.7826346189249
$ws.Range("R4").Value = 522.286516635213
$ws.Range("S4").Value = 21.98354299389757
$ws.Range("T4").Value = -2915.938725396452
$ws.Range("U4").Value = -313.0000453908506
$ws.Range("V4").Value = -245.1738171938969
$ws.Range("W4").Value = -5455.378699940396
$ws.Range("X4").Value = -2745.745190808084
$ws.Range("Y4").Value = -5584.172601264669
$ws.Range("D5").Value = 5373.000859812833
$ws.Range("E5").Value = 3510.376313173212
$ws.Range("F5").Value = 0.6533362649220746
$ws.Range("G5").Value = 1.530605376879352
$ws.Range("H5").Value = 156.2359183937078
$ws.Range("I5").Value = 16.08018431834171
$ws.Range("J5").Value = 13.13678750392228
$ws.Range("K5").Value = 310.6471039720345
$ws.Range("L5").Value = 299.3972496122587
$ws.Range("M5").Value = 187.8434411343187
$ws.Range("N5").Value = 2.784134542991524
$ws.Range("O5").Value = 4.935462251400168
$ws.Range("P5").Value = 0.007548424894594064
$ws.Range("Q5").Value = 348.7826346189249
$ws.Range("R5").Value = 522.286516635213
$ws.Range("S5").Value = 21.98354299389757
$ws.Range("T5").Value = -3121.934233331165
$ws.Range("U5").Value = -316.668224115434
$ws.Range("V5").Value = -262.728201653551
$ws.Range("W5").Value = -5864.159444821766
$ws.Range("X5").Value = -3234.582306051161
$ws.Range("Y5").Value = -5965.961449251277
$ws.Range("D6").Value = 5373.000859812833
$ws.Range("E6").Value = 3752.632387840189
$ws.Range("F6").Value = 0.6984239321284811
$ws.Range("G6").Value = 1.431795151910747
$ws.Range("H6").Value = 166.5306698293716
$ws.Range("I6").Value = 16.26350379249925
$ws.Range("J6").Value = 14.01407859842311
$ws.Range("K6").Value = 331.0761715883855
$ws.Range("L6").Value = 318.477380680386
$ws.Range("M6").Value = 212.273374798242
$ws.Range("N6").Value = 2.784134542991524
$ws.Range("O6").Value = 4.935462251400168
$ws.Range("P6").Value = 0.007548424894594064
$ws.Range("Q6").Value = 348.7826346189249
$ws.Range("R6").Value = 522.286516635213
$ws.Range("S6").Value = 21.98354299389757
$ws.Range("T6").Value = -3327.829262044441
$ws.Range("U6").Value = -320.3346135985848
$ws.Range("V6").Value = -280.2740235435676
$ws.Range("W6").Value = -6272.740797148785
$ws.Range("X6").Value = -3723.180979329627
$ws.Range("Y6").Value = -6347.564070613822
$ws.Range("D7").Value = 5373.000859812833
$ws.Range("E7").Value = 3994.770325143822
$ws.Range("F7").Value = 0.7434896121127696
$ws.Range("G7").Value = 1.34500870450403
$ws.Range("H7").Value = 176.8204009789333
$ws.Range("I7").Value = 16.44673387000694
$ws.Range("J7").Value = 14.89094187761202
$ws.Range("K7").Value = 351.4952768699732
$ws.Range("L7").Value = 337.5482072285376
$ws.Range("M7").Value = 236.6913950848393
$ws.Range("N7").Value = 2.784134542991524
$ws.Range("O7").Value = 4.935462251400168
$ws.Range("P7").Value = 0.007548424894594064
$ws.Range("Q7").Value = 348.7826346189249
$ws.Range("R7").Value = 522.286516635213
$ws.Range("S7").Value = 21.98354299389757
$ws.Range("T7").Value = -3533.623885035675
$ws.Range("U7").Value = -323.9992151487386
$ws.Range("V7").Value = -297.8112891273458
$ws.Range("W7").Value = -6681.122902780538
$ws.Range("X7").Value = -4211.541385061573
$ws.Range("Y7").Value = -6728.980601576855
$ws.Range("D8").Value = 5373.000859812833
$ws.Range("E8").Value = 4236.790211478248
$ws.Range("F8").Value = 0.7885333209542489
$ws.Range("G8").Value = 1.268177226537292
$ws.Range("H8").Value = 187.1051155137393
$ws.Range("I8").Value = 16.6298746162538
$ws.Range("J8").Value = 15.76737765435155
$ws.Range("K8").Value = 371.9044271013699
$ws.Range("L8").Value = 356.6097360621206
$ws.Range("M8").Value = 261.0975107066333
$ws.Range("N8").Value = 2.784134542991524
$ws.Range("O8").Value = 4.935462251400168
$ws.Range("P8").Value = 0.007548424894594064
$ws.Range("Q8").Value = 348.7826346189249
$ws.Range("R8").Value = 522.286516635213
$ws.Range("S8").Value = 21.98354299389757
$ws.Range("T8").Value = -3739.318175731794
$ws.Range("U8").Value = -327.6620300736759
$ws.Range("V8").Value = -315.3400046621364
$ws.Range("W8").Value = -7089.305907408474
$ws.Range("X8").Value = -4699.663697497454
$ws.Range("Y8").Value = -7110.211178248515
